$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Visit date: 13 -> 14 Oct 2021
# ------------------------------------------------------------------
$d.Content.Find.Execute("13 / Oct / 2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "14 / Oct / 2021", 2)

# ------------------------------------------------------------------
# 2. Medicines block: "Metacin ..." -> "Crocin ..." + new "Gelucil ..." line
# ------------------------------------------------------------------
$halfChar = [char]0x00BD

$metacin = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Metacin*") {
        $metacin = $p
        break
    }
}
$metacin.Range.Text = "Crocin  1" + $halfChar + " -- 1 -- 2  for 7 Day(s)"

# re-acquire the paragraph (text replace can invalidate old refs)
$crocin = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Crocin*") {
        $crocin = $p
        break
    }
}

# insert a blank paragraph after it, cloned from a plain (non-list) donor
# paragraph so it does NOT inherit the list numbering / run formatting
$plainDonor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq [char]13) {
        $plainDonor = $p
        break
    }
}
$crocin.Range.InsertParagraphAfter()
$blankAfterCrocin = $crocin.Next()
$blankAfterCrocin.Range.Text = ""
$blankAfterCrocin.Format.Reset()
$blankDonorRange = $plainDonor.Range.Duplicate()
$blankDonorRange.FormattedText.Copy() | Out-Null
$blankAfterCrocin.Range.FormattedText = $blankDonorRange.FormattedText

# insert the new "Gelucil" list item after the blank paragraph, cloned
# from the Crocin list paragraph so list numbering / run formatting match
$blankAfterCrocin.Range.InsertParagraphAfter()
$gelucil = $blankAfterCrocin.Next()
$gelucil.Range.Text = "Gelucil  2" + $halfChar + " -- 0 -- 3  for 2 Weeks(s)"

# ------------------------------------------------------------------
# 3. Advice block: "Avoid cold drink" -> "Avoid Sugar ..." + 4 new lines
# ------------------------------------------------------------------
$avoid = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Avoid cold drink*") {
        $avoid = $p
        break
    }
}
$avoid.Range.Text = "Avoid Sugar / Gur / Honey Fruit Juice and Soft Drink"

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Avoid Sugar*") {
        $anchor = $p
        break
    }
}

$newItems = @(
    "Any fruit only 100g per day",
    "Balanced Diet",
    "Take steam 2 / 3 time a day",
    "Physiotherapy"
)
foreach ($item in $newItems) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()
    $newPara.Range.Text = $item
    $anchor = $newPara
}

# ------------------------------------------------------------------
# 4. Test block: collapse "Hb1aC after 3 months" / "ECG" / "xray"
#    into a single "Rem 1" line
# ------------------------------------------------------------------
$hb1ac = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Hb1aC*") {
        $hb1ac = $p
        break
    }
}
$hb1ac.Range.Text = "Rem 1"

$toRemove = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -like "ECG*") -or ($t -like "xray*")) {
        $toRemove += $p
    }
}
for ($i = $toRemove.Count - 1; $i -ge 0; $i--) {
    $toRemove[$i].Range.Delete()
}

# ------------------------------------------------------------------
# 5. Next review date: 27 -> 28 / 10 / 2021
# ------------------------------------------------------------------
$d.Content.Find.Execute("27 / 10 / 2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "28 / 10 / 2021", 2)
